# Logbook update: add 4 new rows (29-32) of log entries.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "MM/DD/YY"

# Row 29 - 11/29/2019 - pipeline/ROS integration + client service fix
$ws.Cells.Item(29, 1).Value = 43798
$ws.Cells.Item(29, 1).NumberFormat = $dateFormat
$ws.Cells.Item(29, 2).Value = "Integrated pipeline into ROS"
$ws.Cells.Item(29, 3).Value = "Fixed client service with images"

# Row 30 - 12/2/2019 - client service integrated into main (column B left empty)
$ws.Cells.Item(30, 1).Value = 43801
$ws.Cells.Item(30, 1).NumberFormat = $dateFormat
$ws.Cells.Item(30, 3).Value = "Integrated client service into main"

# Row 31 - 12/3/2019 - date only
$ws.Cells.Item(31, 1).Value = 43802
$ws.Cells.Item(31, 1).NumberFormat = $dateFormat

# Row 32 - 12/4/2019 - date only
$ws.Cells.Item(32, 1).Value = 43803
$ws.Cells.Item(32, 1).NumberFormat = $dateFormat

# Update the selected cell to match the author's final cursor position
$ws.Range("B31").Select()
